$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 was storing numeric-looking values as text; convert them to real numbers.
$ws.Range("B6").Value = 4141
$ws.Range("C6").Value = 21881
$ws.Range("D6").Value = 4142
$ws.Range("E6").Value = 22183
$ws.Range("F6").Value = 302
$ws.Range("G6").Value = 101.38

# Insert a new data row (row 7) below, keeping the same text formatting
# pattern as the rest of the list (values stored as text, e.g. "4000.00").
$ws.Range("A7:G7").NumberFormat = "@"
$ws.Range("A7").Value = "03/08/2023"
$ws.Range("B7").Value = "4000.00"
$ws.Range("C7").Value = "25881.00"
$ws.Range("D7").Value = "4000.00"
$ws.Range("E7").Value = "26183.00"
$ws.Range("F7").Value = "302.00"
$ws.Range("G7").Value = "101.17"
$ws.Range("A7:G7").Style = "Normal"
